# Updated symbol list on Wed Feb  8 22:52:02 UTC 2023 with GitHub Actions
#
# Refreshes the coin-ranking snapshot on Sheet1 (rows 2-51): most rows pick
# up a new Price (col D) and/or Volume(1h) (col E) reading, and the
# FTXToken..LEO block (rows 6-17) shifts down one slot to make room for a
# new top entry (GateToken) together with its Link (col C) and Price.
#
# Every data cell in columns B-G on this sheet is stored as literal text
# (e.g. column D holds strings like "327.98" or "0.0001300", column E holds
# strings like "-1.27%"), never as a real number/percentage. A plain
# `Range.Value = "..."` assignment lets Excel's type-inference kick in and
# silently turn such strings into numbers (or percentage-formatted
# numbers), changing both the stored type and the cell style. To avoid
# that, numeric-looking columns (D and E) are written with the cell's
# NumberFormat temporarily forced to "@" (Text); the cell style is then
# reset to "Normal" afterwards so no residual formatting is left behind.
# Columns B and C (coin name / link) are plain text already, so they are
# assigned directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "327.98" },
    @{ Cell = "E2"; Value = "-1.27%" },
    @{ Cell = "D3"; Value = "44.33" },
    @{ Cell = "E3"; Value = "-0.95%" },
    @{ Cell = "D4"; Value = "5.405" },
    @{ Cell = "E4"; Value = "-2.58%" },
    @{ Cell = "E5"; Value = "1.10%" },
    @{ Cell = "B6"; Value = "GateToken" },
    @{ Cell = "C6"; Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt" },
    @{ Cell = "D6"; Value = "4.425" },
    @{ Cell = "E6"; Value = "-0.12%" },
    @{ Cell = "B7"; Value = "FTXToken" },
    @{ Cell = "C7"; Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt" },
    @{ Cell = "D7"; Value = "1.936" },
    @{ Cell = "E7"; Value = "-4.90%" },
    @{ Cell = "B8"; Value = "MXToken" },
    @{ Cell = "C8"; Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx" },
    @{ Cell = "D8"; Value = "0.9726" },
    @{ Cell = "E8"; Value = "-0.71%" },
    @{ Cell = "B9"; Value = "BTSEToken" },
    @{ Cell = "C9"; Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse" },
    @{ Cell = "D9"; Value = "2.510" },
    @{ Cell = "E9"; Value = "-4.41%" },
    @{ Cell = "B10"; Value = "LiechtensteinCryptoassetsExchange" },
    @{ Cell = "C10"; Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx" },
    @{ Cell = "D10"; Value = "0.1135" },
    @{ Cell = "E10"; Value = "0.78%" },
    @{ Cell = "B11"; Value = "WazirX" },
    @{ Cell = "C11"; Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx" },
    @{ Cell = "D11"; Value = "0.1901" },
    @{ Cell = "E11"; Value = "-0.43%" },
    @{ Cell = "B12"; Value = "MandalaExchangeToken" },
    @{ Cell = "C12"; Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx" },
    @{ Cell = "D12"; Value = "0.09709" },
    @{ Cell = "E12"; Value = "-3.43%" },
    @{ Cell = "B13"; Value = "BitrueCoin" },
    @{ Cell = "C13"; Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr" },
    @{ Cell = "D13"; Value = "0.04611" },
    @{ Cell = "E13"; Value = "-1.49%" },
    @{ Cell = "B14"; Value = "BitMartToken" },
    @{ Cell = "C14"; Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx" },
    @{ Cell = "D14"; Value = "0.1062" },
    @{ Cell = "E14"; Value = "0.36%" },
    @{ Cell = "B15"; Value = "BitForexToken" },
    @{ Cell = "C15"; Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf" },
    @{ Cell = "D15"; Value = "0.001293" },
    @{ Cell = "E15"; Value = "1.92%" },
    @{ Cell = "B16"; Value = "TigerCash" },
    @{ Cell = "C16"; Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch" },
    @{ Cell = "D16"; Value = "0.005965" },
    @{ Cell = "E16"; Value = "-1.97%" },
    @{ Cell = "B17"; Value = "LEO" },
    @{ Cell = "C17"; Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo" },
    @{ Cell = "D17"; Value = "3.361" },
    @{ Cell = "E17"; Value = "0.02%" },
    @{ Cell = "D18"; Value = "0.3358" },
    @{ Cell = "E18"; Value = "0.13%" },
    @{ Cell = "D19"; Value = "8.312" },
    @{ Cell = "E19"; Value = "-19.08%" },
    @{ Cell = "E20"; Value = "0.39%" },
    @{ Cell = "E21"; Value = "6.75%" },
    @{ Cell = "D22"; Value = "0.04175" },
    @{ Cell = "E22"; Value = "1.56%" },
    @{ Cell = "E23"; Value = "-4.77%" },
    @{ Cell = "D24"; Value = "0.004449" },
    @{ Cell = "E24"; Value = "1.15%" },
    @{ Cell = "D25"; Value = "0.0001299" },
    @{ Cell = "E25"; Value = "1.55%" },
    @{ Cell = "D26"; Value = "0.0002979" },
    @{ Cell = "E26"; Value = "-20.22%" },
    @{ Cell = "D38"; Value = "0.02712" },
    @{ Cell = "E38"; Value = "-3.12%" },
    @{ Cell = "D39"; Value = "0.05632" },
    @{ Cell = "E39"; Value = "-1.99%" },
    @{ Cell = "D40"; Value = "0.007865" },
    @{ Cell = "E40"; Value = "3.20%" },
    @{ Cell = "D41"; Value = "0.1413" },
    @{ Cell = "E41"; Value = "-0.75%" },
    @{ Cell = "D42"; Value = "0.007293" },
    @{ Cell = "E42"; Value = "-3.34%" },
    @{ Cell = "D43"; Value = "0.002039" },
    @{ Cell = "E43"; Value = "3.40%" },
    @{ Cell = "D44"; Value = "0.008711" },
    @{ Cell = "E44"; Value = "8.50%" },
    @{ Cell = "D45"; Value = "0.3504" },
    @{ Cell = "D46"; Value = "0.00006903" },
    @{ Cell = "E46"; Value = "-1.56%" },
    @{ Cell = "E47"; Value = "0.17%" },
    @{ Cell = "D48"; Value = "0.003485" },
    @{ Cell = "E48"; Value = "-2.56%" },
    @{ Cell = "D49"; Value = "0.003530" },
    @{ Cell = "E49"; Value = "39.87%" },
    @{ Cell = "D50"; Value = "0.00002100" },
    @{ Cell = "E50"; Value = "0.17%" },
    @{ Cell = "D51"; Value = "0.0002000" },
    @{ Cell = "E51"; Value = "0.17%" }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $col = $u.Cell -replace '[0-9]+$', ''

    if ($col -eq "D" -or $col -eq "E") {
        # Numeric/percent-looking text: force Text format so Excel keeps
        # the literal string instead of coercing it to a number, then
        # restore the default "Normal" style so no formatting sticks.
        $cell.NumberFormat = "@"
        $cell.Value = $u.Value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $u.Value
    }
}
